# prices and index update
# The "Table 0" sheet holds the result of a Power-Query refresh (montenegro
# fuel prices). The refreshed query dropped the oldest 10 days
# (26.01.2025-04.02.2025) and appended 10 new days (26.02.2025-07.03.2025),
# while a few of the rolling price columns shifted too. Reproduce the
# refreshed table contents directly on the worksheet cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table 0")

$dates = @("05.02.2025","06.02.2025","07.02.2025","08.02.2025","09.02.2025","10.02.2025","11.02.2025","12.02.2025","13.02.2025","14.02.2025","15.02.2025","16.02.2025","17.02.2025","18.02.2025","19.02.2025","20.02.2025","21.02.2025","22.02.2025","23.02.2025","24.02.2025","25.02.2025","26.02.2025","27.02.2025","28.02.2025","01.03.2025","02.03.2025","03.03.2025","04.03.2025","05.03.2025","06.03.2025","07.03.2025")

$gasoline = @("2,93","2,93","2,93","2,93","2,93","2,93","2,93","2,97","2,97","2,97","2,97","2,97","2,97","2,97","2,97","2,97","2,97","2,97","2,97","2,97","2,97","2,97","2,97","2,97","2,97","2,97","2,97","2,97","2,97","2,97","2,97")

$diesel = @("2,80","2,80","2,80","2,80","2,80","2,80","2,80","2,80","2,80","2,80","2,80","2,80","2,80","2,80","2,80","2,80","2,80","2,80","2,80","2,80","2,80","2,82","2,82","2,82","2,82","2,82","2,82","2,82","2,82","2,82","2,82")

$gas = @("1,27","1,27","1,27","1,27","1,27","1,27","1,27","1,27","1,27","1,27","1,27","1,27","1,27","1,27","1,27","1,27","1,27","1,27","1,27","1,27","1,27","1,27","1,27","1,27","1,27","1,27","1,27","1,27","1,27","1,27","1,27")

# Non-breaking space: the source data always trails the numeric price
# strings with U+00A0 so Excel keeps them as text instead of parsing them
# as (locale-comma) numbers.
$nbsp = [char]0x00A0

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    # Leading apostrophe forces the date-looking string to be stored as
    # literal text (matching the existing column, which is plain text,
    # not a real date serial).
    $ws.Cells.Item($row, 1).Value = "'" + $dates[$i]
    $ws.Cells.Item($row, 2).Value = $gasoline[$i] + $nbsp
    $ws.Cells.Item($row, 3).Value = $diesel[$i] + $nbsp
    $ws.Cells.Item($row, 4).Value = $gas[$i] + $nbsp
}
